$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
  "04-08-2021",
  "05-08-2021",
  "06-08-2021",
  "09-08-2021",
  "10-08-2021",
  "11-08-2021",
  "12-08-2021",
  "13-08-2021",
  "16-08-2021",
  "17-08-2021",
  "18-08-2021",
  "19-08-2021",
  "20-08-2021",
  "23-08-2021",
  "24-08-2021",
  "25-08-2021",
  "26-08-2021",
  "27-08-2021",
  "30-08-2021",
  "31-08-2021",
  "01-09-2021",
  "02-09-2021",
  "03-09-2021"
)

$data = @(
  @(0.67,0.95,1.37,1.95,-0.44),
  @(0.67,0.99,1.42,1.96,-0.47),
  @(0.71,1.07,1.47,1.99,-0.53),
  @(0.73,1.07,1.48,2,-0.55),
  @(0.72,1.09,1.52,2.03,-0.55),
  @(0.72,1.1,1.57,2.1,-0.5),
  @(0.72,1.13,1.6,2.14,-0.55),
  @(0.72,1.13,1.6,2.14,-0.5),
  @(0.72,1.13,1.6,2.14,-0.6899999999999999),
  @(0.72,1.12,1.6,2.16,-0.72),
  @(0.72,1.17,1.61,2.29,-0.6899999999999999),
  @(0.74,1.24,1.64,2.34,-0.65),
  @(0.74,1.26,1.71,2.32,-0.64),
  @(0.75,1.27,1.72,2.34,-0.62),
  @(0.77,1.28,1.73,2.33,-0.62),
  @(0.78,1.28,1.73,2.34,-0.64),
  @(0.84,1.31,1.76,2.34,-0.65),
  @(0.86,1.32,1.77,2.34,-0.65),
  @(0.86,1.29,1.77,2.38,-0.66),
  @(0.88,1.27,1.77,2.35,-0.66),
  @(1.39,1.64,2.07,2.61,-0.57),
  @(1.45,1.84,2.22,2.9,-0.47),
  @(1.45,1.87,2.27,2.93,-0.44)
)

$startRow = 149
for ($i = 0; $i -lt $dates.Length; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $dates[$i]
  $ws.Cells.Item($r, 2).Value = $data[$i][0]
  $ws.Cells.Item($r, 3).Value = $data[$i][1]
  $ws.Cells.Item($r, 4).Value = $data[$i][2]
  $ws.Cells.Item($r, 5).Value = $data[$i][3]
  $ws.Cells.Item($r, 6).Value = $data[$i][4]
}
